# TC01_Trials_Filter_Gender-Male.xlsx — add the Neo4j Cypher query used by
# the WebData/Neo4jData comparison, resize the query cell's row to fit the
# wrapped text, and move the sheet's active selection to the new cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the (until now empty) query cell for the "dbExcel" / "query" row.
# Populate it with the Cypher query text; the cell already carries the
# wrap-text style (s="1") from the template.
$ws.Range("A2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.gender IN ['MALE'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

# Grow row 2 so the long, wrapped query text is fully visible.
$ws.Rows.Item(2).RowHeight = 87

# The author's saved selection moved from B7 to the newly-filled A2, and the
# view no longer needs to be scrolled to keep a far-right cell (B7) in frame.
[void]$ws.Range("A2").Select()
